$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = [double]"0.003451937919490256"
$ws.Range("C3").Value = [double]"0.002174424256882952"
$ws.Range("C4").Value = [double]"0.002063949932462754"
$ws.Range("C5").Value = [double]"0.001278763926158885"
$ws.Range("C6").Value = [double]"0.0005265073138351496"
$ws.Range("C7").Value = [double]"0.0005029599529687292"
$ws.Range("C8").Value = [double]"0.0005052532361994914"
$ws.Range("C9").Value = [double]"0.000423126794535495"
$ws.Range("C10").Value = [double]"0.0003997944618417414"
$ws.Range("C11").Value = [double]"0.0003097026208987869"
$ws.Range("C12").Value = [double]"0.0002668173377425682"
$ws.Range("C13").Value = [double]"0.0002367816917222445"
$ws.Range("C14").Value = [double]"0.0002184855760046509"
$ws.Range("C15").Value = [double]"0.0002216883783231014"
$ws.Range("C16").Value = [double]"0.0002041136020256109"
$ws.Range("C17").Value = [double]"0.0001625398327072145"
$ws.Range("C18").Value = [double]"0.0001416810574892402"
$ws.Range("C19").Value = [double]"0.000132359584053288"
$ws.Range("C20").Value = [double]"0.0001008854184642631"
$ws.Range("C21").Value = [double]"0.0001024871207545312"
$ws.Range("C22").Value = [double]"8.958136837734933e-05"
$ws.Range("C23").Value = [double]"8.615430263815464e-05"
$ws.Range("C24").Value = [double]"9.085692717647688e-05"
$ws.Range("C25").Value = [double]"8.261384398590614e-05"
$ws.Range("C26").Value = [double]"7.676615832362616e-05"
$ws.Range("C27").Value = [double]"8.173308377374262e-05"
$ws.Range("C28").Value = [double]"7.657206567481889e-05"
$ws.Range("C29").Value = [double]"7.198766389611014e-05"
$ws.Range("C30").Value = [double]"6.758690755563884e-05"
$ws.Range("C31").Value = [double]"5.442784661287518e-05"
$ws.Range("C32").Value = [double]"5.297811058123772e-05"
$ws.Range("C33").Value = [double]"4.527378880840322e-05"
$ws.Range("C34").Value = [double]"4.243899784766316e-05"
$ws.Range("C35").Value = [double]"4.154300782482745e-05"
$ws.Range("C36").Value = [double]"4.169304120285888e-05"
$ws.Range("C37").Value = [double]"4.333149484771969e-05"
$ws.Range("C38").Value = [double]"3.415737606811432e-05"
$ws.Range("C39").Value = [double]"3.365592709574874e-05"
$ws.Range("C40").Value = [double]"2.97136497841402e-05"
$ws.Range("C41").Value = [double]"3.06058726696013e-05"
$ws.Range("C42").Value = [double]"2.764108534118914e-05"
$ws.Range("C43").Value = [double]"2.599747166957288e-05"
$ws.Range("C44").Value = [double]"1.882879100141476e-05"
$ws.Range("C45").Value = [double]"1.927092727293902e-05"
$ws.Range("C46").Value = [double]"1.827180048370134e-05"
$ws.Range("C47").Value = [double]"1.890231684228599e-05"
$ws.Range("C48").Value = [double]"1.899634361467679e-05"
$ws.Range("C49").Value = [double]"1.875840144161789e-05"
$ws.Range("C50").Value = [double]"1.777120928459473e-05"
$ws.Range("C51").Value = [double]"1.82408418420414e-05"
$ws.Range("C52").Value = [double]"1.688504416246849e-05"
$ws.Range("C53").Value = [double]"1.771945217078104e-05"
$ws.Range("C54").Value = [double]"1.543440494995655e-05"
$ws.Range("C55").Value = [double]"1.637921612600122e-05"
$ws.Range("C56").Value = [double]"1.430661841067594e-05"
$ws.Range("C57").Value = [double]"1.476946652795001e-05"
$ws.Range("C58").Value = [double]"1.486528247531807e-05"
$ws.Range("C59").Value = [double]"1.46860767990808e-05"
$ws.Range("C60").Value = [double]"1.389855959298695e-05"
$ws.Range("C61").Value = [double]"1.387583807824674e-05"
$ws.Range("C62").Value = [double]"1.363481643619575e-05"
$ws.Range("C63").Value = [double]"1.289061536843044e-05"
$ws.Range("C64").Value = [double]"1.238307962696959e-05"
$ws.Range("C65").Value = [double]"1.216328619537913e-05"
$ws.Range("C66").Value = [double]"1.089681647703911e-05"
$ws.Range("C67").Value = [double]"1.044997807021348e-05"
$ws.Range("C68").Value = [double]"9.546882602231099e-06"
$ws.Range("C69").Value = [double]"9.092133689532084e-06"
$ws.Range("C70").Value = [double]"8.724301742754743e-06"
$ws.Range("C71").Value = [double]"7.102920859533855e-06"
$ws.Range("C72").Value = [double]"6.969569692280081e-06"
$ws.Range("C73").Value = [double]"5.242882679187686e-06"
$ws.Range("C74").Value = [double]"4.66194168901706e-06"
$ws.Range("C75").Value = [double]"4.738607483292845e-06"
$ws.Range("C76").Value = [double]"3.966300150926133e-06"
$ws.Range("C77").Value = [double]"4.032357452569731e-06"
$ws.Range("C78").Value = [double]"3.641995801931504e-06"
$ws.Range("C79").Value = [double]"3.628809335402092e-06"
$ws.Range("C80").Value = [double]"3.441908802092485e-06"
$ws.Range("C81").Value = [double]"3.068620634200842e-06"
$ws.Range("C82").Value = [double]"2.348573839717578e-06"
$ws.Range("C83").Value = [double]"2.269776099189785e-06"
$ws.Range("C84").Value = [double]"2.177048988093499e-06"
$ws.Range("C85").Value = [double]"1.985663425270344e-06"
$ws.Range("C86").Value = [double]"1.666568415716128e-06"
$ws.Range("C87").Value = [double]"1.533878974283268e-06"
$ws.Range("C88").Value = [double]"1.323603591546073e-06"
$ws.Range("C89").Value = [double]"1.253227205356783e-06"
$ws.Range("C90").Value = [double]"3.736497162409213e-07"
$ws.Range("C91").Value = [double]"5.512647892298724e-07"
$ws.Range("C92").Value = [double]"2.301740199921175e-08"
$ws.Range("C93").Value = [double]"1.8729973631225e-08"
$ws.Range("C94").Value = [double]"1.351197923590519e-08"
$ws.Range("C95").Value = [double]"1.675691880274835e-08"
$ws.Range("C96").Value = [double]"4.863780410225193e-09"
$ws.Range("C97").Value = [double]"3.274894539521147e-09"
$ws.Range("C98").Value = [double]"1.236383399830115e-08"
$ws.Range("C99").Value = [double]"5.807800180070509e-09"
$ws.Range("C100").Value = [double]"7.746352843261796e-09"
$ws.Range("C101").Value = [double]"6.060838160530408e-09"
$ws.Range("C102").Value = [double]"1.581291491316772e-08"
$ws.Range("C103").Value = [double]"4.190068683019831e-10"

$dataRange = $ws.Range("A2:C103")
$keyRange = $ws.Range("C2:C103")
$dataRange.Sort($keyRange, 2)
